$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178:240 down to 179:241
$ws.Rows.Item(178).Insert()

# Populate the new row 178 - same categorical/reference data as the row below it
# (Terminal Hortofrutícola Agro Chillán / Ñuble / Piña / Caramelo / Segunda / Ecuador),
# with the new date, volume and price figures from this week's update.
$ws.Cells.Item(178, 1).Value = 7
$ws.Cells.Item(178, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(178, 3).Value = "Ñuble"
$ws.Cells.Item(178, 4).Value = 44809
$ws.Cells.Item(178, 5).Value = 16
$ws.Cells.Item(178, 6).Value = "Fruta"
$ws.Cells.Item(178, 7).Value = 100108
$ws.Cells.Item(178, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(178, 9).Value = 100108005
$ws.Cells.Item(178, 10).Value = "Piña"
$ws.Cells.Item(178, 11).Value = "Caramelo"
$ws.Cells.Item(178, 12).Value = "Segunda"
$ws.Cells.Item(178, 13).Value = 60
$ws.Cells.Item(178, 14).Value = 19000
$ws.Cells.Item(178, 15).Value = 19000
$ws.Cells.Item(178, 16).Value = 19000
$ws.Cells.Item(178, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(178, 18).Value = "Ecuador"
$ws.Cells.Item(178, 19).Value = 1357
$ws.Cells.Item(178, 20).Value = 14
